# This workbook ("Hortaliza, Macroferia Regional de Talca - Choclo") keeps
# one price-report row per line, always shifting older rows down as a new
# weekly report is added at the top of the data block (row 55, right after
# the row that is currently the most recent at row 54).
#
# The edit: insert one new row at row 55 (pushing the existing rows 55-149
# down to 56-150, which grows the sheet from A1:R149 to A1:R150), and
# populate the newly inserted row 55 with this week's report data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 55:149 down to 56:150, leaving row 55 empty (except it
# inherits the date column's number format from the surrounding cells).
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with the new weekly record.
$ws.Cells.Item(55, 1).Value = 5
$ws.Cells.Item(55, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(55, 3).Value = "Maule"
$ws.Cells.Item(55, 4).Value = 44557
$ws.Cells.Item(55, 5).Value = 7
$ws.Cells.Item(55, 6).Value = 100112024
$ws.Cells.Item(55, 7).Value = "Choclo"
$ws.Cells.Item(55, 8).Value = "Choclero"
$ws.Cells.Item(55, 9).Value = "Segunda"
$ws.Cells.Item(55, 10).Value = 60000
$ws.Cells.Item(55, 11).Value = 200
$ws.Cells.Item(55, 12).Value = 200
$ws.Cells.Item(55, 13).Value = 200
$ws.Cells.Item(55, 14).Value = "$/unidad"
$ws.Cells.Item(55, 15).Value = "Región del Maule"
$ws.Cells.Item(55, 16).Value = 200
$ws.Cells.Item(55, 17).Value = 1
$ws.Cells.Item(55, 18).Value = "Hortaliza"
